$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 19:52"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 829392
$ws.Range("C4").Value = 10648
$ws.Range("E4").Value = 699823
$ws.Range("G4").Value = 831
$ws.Range("H4").Value = 46149

# Canada (row 16)
$ws.Range("B16").Value = 39805
$ws.Range("C16").Value = 1383
$ws.Range("E16").Value = 24192
$ws.Range("G16").Value = 132
$ws.Range("H16").Value = 1966

# Irlanda (row 22)
$ws.Range("B22").Value = 16671
$ws.Range("C22").Value = 631
$ws.Range("E22").Value = 6669
$ws.Range("F22").Value = 147
$ws.Range("G22").Value = 39
$ws.Range("H22").Value = 769

# Rumania (row 34)
$ws.Range("E34").Value = 6785
$ws.Range("G34").Value = 21
$ws.Range("H34").Value = 519

# Luxemburgo (row 53)
$ws.Range("B53").Value = 3654
$ws.Range("C53").Value = 36
$ws.Range("D53").Value = 711
$ws.Range("E53").Value = 2863
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 80

# Jordania (row 108)
$ws.Range("B108").Value = 435
$ws.Range("C108").Value = 7
$ws.Range("E108").Value = 131

# Swap Mayotte / Sri Lanka ordering (rows 114 & 115) and update Sri Lanka's
# figures. Row 114 becomes Sri Lanka (with updated counts); row 115 becomes
# Mayotte (keeping its previous counts).
$ws.Range("A114").Value = "Sri Lanka"
$ws.Range("B114").Value = 328
$ws.Range("C114").Value = 18
$ws.Range("D114").Value = 105
$ws.Range("E114").Value = 216
$ws.Range("F114").Value = 2
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 7

$ws.Range("A115").Value = "Mayotte"
$ws.Range("B115").Value = 326
$ws.Range("C115").Value = 15
$ws.Range("D115").Value = 125
$ws.Range("E115").Value = 197
$ws.Range("F115").Value = 4
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 4

# Santa Lucia (row 190)
$ws.Range("D190").Value = 15
$ws.Range("E190").Value = 0
